# Updated symbol list: refreshed prices, volume(1h) percentages and the
# "Hora" (hour) column for every coin row, plus a rotation of the
# Coin/Link pair for rows 9-15 (each row takes the coin that was one row
# below it, with row 15 receiving what row 9 used to hold).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry describes one cell that changes. "IsNumericText" marks cells
# whose new text looks like a number/percentage (e.g. "256.64", "0.49%",
# "13") -- for those we briefly force a text NumberFormat so Excel stores
# the literal string instead of silently converting it to a numeric
# value, then clear the format again so no stray style is left behind.
$changes = @(
    @{ Cell = "D2"; Value = "256.64"; IsNumericText = $true },
    @{ Cell = "E2"; Value = "0.49%"; IsNumericText = $true },
    @{ Cell = "G2"; Value = "13"; IsNumericText = $true },
    @{ Cell = "D3"; Value = "27.14"; IsNumericText = $true },
    @{ Cell = "E3"; Value = "-2.23%"; IsNumericText = $true },
    @{ Cell = "G3"; Value = "13"; IsNumericText = $true },
    @{ Cell = "D4"; Value = "4.712"; IsNumericText = $true },
    @{ Cell = "E4"; Value = "-9.91%"; IsNumericText = $true },
    @{ Cell = "G4"; Value = "13"; IsNumericText = $true },
    @{ Cell = "D5"; Value = "0.05922"; IsNumericText = $true },
    @{ Cell = "E5"; Value = "0.86%"; IsNumericText = $true },
    @{ Cell = "G5"; Value = "13"; IsNumericText = $true },
    @{ Cell = "E6"; Value = "-0.67%"; IsNumericText = $true },
    @{ Cell = "G6"; Value = "13"; IsNumericText = $true },
    @{ Cell = "D7"; Value = "0.8684"; IsNumericText = $true },
    @{ Cell = "E7"; Value = "0.03%"; IsNumericText = $true },
    @{ Cell = "G7"; Value = "13"; IsNumericText = $true },
    @{ Cell = "D8"; Value = "0.9445"; IsNumericText = $true },
    @{ Cell = "E8"; Value = "-0.30%"; IsNumericText = $true },
    @{ Cell = "G8"; Value = "13"; IsNumericText = $true },
    @{ Cell = "B9"; Value = "One"; IsNumericText = $false },
    @{ Cell = "C9"; Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"; IsNumericText = $false },
    @{ Cell = "D9"; Value = "0.0006029"; IsNumericText = $true },
    @{ Cell = "E9"; Value = "-0.34%"; IsNumericText = $true },
    @{ Cell = "G9"; Value = "13"; IsNumericText = $true },
    @{ Cell = "B10"; Value = "WazirX"; IsNumericText = $false },
    @{ Cell = "C10"; Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"; IsNumericText = $false },
    @{ Cell = "D10"; Value = "0.1409"; IsNumericText = $true },
    @{ Cell = "E10"; Value = "-0.08%"; IsNumericText = $true },
    @{ Cell = "G10"; Value = "13"; IsNumericText = $true },
    @{ Cell = "B11"; Value = "LiechtensteinCryptoassetsExchange"; IsNumericText = $false },
    @{ Cell = "C11"; Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"; IsNumericText = $false },
    @{ Cell = "D11"; Value = "0.03885"; IsNumericText = $true },
    @{ Cell = "E11"; Value = "11.71%"; IsNumericText = $true },
    @{ Cell = "G11"; Value = "13"; IsNumericText = $true },
    @{ Cell = "B12"; Value = "MandalaExchangeToken"; IsNumericText = $false },
    @{ Cell = "C12"; Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"; IsNumericText = $false },
    @{ Cell = "D12"; Value = "0.07098"; IsNumericText = $true },
    @{ Cell = "E12"; Value = "-0.94%"; IsNumericText = $true },
    @{ Cell = "G12"; Value = "13"; IsNumericText = $true },
    @{ Cell = "B13"; Value = "BitrueCoin"; IsNumericText = $false },
    @{ Cell = "C13"; Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"; IsNumericText = $false },
    @{ Cell = "D13"; Value = "0.03213"; IsNumericText = $true },
    @{ Cell = "E13"; Value = "0.86%"; IsNumericText = $true },
    @{ Cell = "G13"; Value = "13"; IsNumericText = $true },
    @{ Cell = "B14"; Value = "BitMartToken"; IsNumericText = $false },
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"; IsNumericText = $false },
    @{ Cell = "D14"; Value = "0.09255"; IsNumericText = $true },
    @{ Cell = "E14"; Value = "0.40%"; IsNumericText = $true },
    @{ Cell = "G14"; Value = "13"; IsNumericText = $true },
    @{ Cell = "B15"; Value = "BitForexToken"; IsNumericText = $false },
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"; IsNumericText = $false },
    @{ Cell = "D15"; Value = "0.001548"; IsNumericText = $true },
    @{ Cell = "E15"; Value = "-0.84%"; IsNumericText = $true },
    @{ Cell = "G15"; Value = "13"; IsNumericText = $true },
    @{ Cell = "D16"; Value = "0.006048"; IsNumericText = $true },
    @{ Cell = "E16"; Value = "3.34%"; IsNumericText = $true },
    @{ Cell = "G16"; Value = "13"; IsNumericText = $true },
    @{ Cell = "D17"; Value = "3.513"; IsNumericText = $true },
    @{ Cell = "E17"; Value = "0.40%"; IsNumericText = $true },
    @{ Cell = "G17"; Value = "13"; IsNumericText = $true },
    @{ Cell = "D18"; Value = "3.196"; IsNumericText = $true },
    @{ Cell = "E18"; Value = "-0.99%"; IsNumericText = $true },
    @{ Cell = "G18"; Value = "13"; IsNumericText = $true },
    @{ Cell = "E19"; Value = "0.64%"; IsNumericText = $true },
    @{ Cell = "G19"; Value = "13"; IsNumericText = $true },
    @{ Cell = "D20"; Value = "0.3141"; IsNumericText = $true },
    @{ Cell = "E20"; Value = "-1.11%"; IsNumericText = $true },
    @{ Cell = "G20"; Value = "13"; IsNumericText = $true },
    @{ Cell = "E21"; Value = "-0.76%"; IsNumericText = $true },
    @{ Cell = "G21"; Value = "13"; IsNumericText = $true },
    @{ Cell = "D22"; Value = "3.824"; IsNumericText = $true },
    @{ Cell = "E22"; Value = "8.52%"; IsNumericText = $true },
    @{ Cell = "G22"; Value = "13"; IsNumericText = $true },
    @{ Cell = "D23"; Value = "0.04219"; IsNumericText = $true },
    @{ Cell = "E23"; Value = "1.24%"; IsNumericText = $true },
    @{ Cell = "G23"; Value = "13"; IsNumericText = $true },
    @{ Cell = "E24"; Value = "0.38%"; IsNumericText = $true },
    @{ Cell = "G24"; Value = "13"; IsNumericText = $true },
    @{ Cell = "E25"; Value = "-0.47%"; IsNumericText = $true },
    @{ Cell = "G25"; Value = "13"; IsNumericText = $true },
    @{ Cell = "D26"; Value = "0.004292"; IsNumericText = $true },
    @{ Cell = "E26"; Value = "-10.51%"; IsNumericText = $true },
    @{ Cell = "G26"; Value = "13"; IsNumericText = $true },
    @{ Cell = "D27"; Value = "0.0001200"; IsNumericText = $true },
    @{ Cell = "E27"; Value = "-0.03%"; IsNumericText = $true },
    @{ Cell = "G27"; Value = "13"; IsNumericText = $true },
    @{ Cell = "E28"; Value = "2.39%"; IsNumericText = $true },
    @{ Cell = "G28"; Value = "13"; IsNumericText = $true },
    @{ Cell = "G29"; Value = "13"; IsNumericText = $true },
    @{ Cell = "G30"; Value = "13"; IsNumericText = $true },
    @{ Cell = "G31"; Value = "13"; IsNumericText = $true },
    @{ Cell = "G32"; Value = "13"; IsNumericText = $true },
    @{ Cell = "G33"; Value = "13"; IsNumericText = $true },
    @{ Cell = "G34"; Value = "13"; IsNumericText = $true },
    @{ Cell = "G35"; Value = "13"; IsNumericText = $true },
    @{ Cell = "G36"; Value = "13"; IsNumericText = $true },
    @{ Cell = "G37"; Value = "13"; IsNumericText = $true },
    @{ Cell = "G38"; Value = "13"; IsNumericText = $true },
    @{ Cell = "G39"; Value = "13"; IsNumericText = $true },
    @{ Cell = "D40"; Value = "0.03828"; IsNumericText = $true },
    @{ Cell = "E40"; Value = "0.59%"; IsNumericText = $true },
    @{ Cell = "G40"; Value = "13"; IsNumericText = $true },
    @{ Cell = "D41"; Value = "0.006249"; IsNumericText = $true },
    @{ Cell = "E41"; Value = "64.67%"; IsNumericText = $true },
    @{ Cell = "G41"; Value = "13"; IsNumericText = $true },
    @{ Cell = "D42"; Value = "0.1102"; IsNumericText = $true },
    @{ Cell = "E42"; Value = "-0.10%"; IsNumericText = $true },
    @{ Cell = "G42"; Value = "13"; IsNumericText = $true },
    @{ Cell = "D43"; Value = "0.001900"; IsNumericText = $true },
    @{ Cell = "E43"; Value = "-20.28%"; IsNumericText = $true },
    @{ Cell = "G43"; Value = "13"; IsNumericText = $true },
    @{ Cell = "D44"; Value = "0.01149"; IsNumericText = $true },
    @{ Cell = "E44"; Value = "17.42%"; IsNumericText = $true },
    @{ Cell = "G44"; Value = "13"; IsNumericText = $true },
    @{ Cell = "E45"; Value = "3.01%"; IsNumericText = $true },
    @{ Cell = "G45"; Value = "13"; IsNumericText = $true },
    @{ Cell = "E46"; Value = "-0.01%"; IsNumericText = $true },
    @{ Cell = "G46"; Value = "13"; IsNumericText = $true },
    @{ Cell = "D47"; Value = "0.08049"; IsNumericText = $true },
    @{ Cell = "E47"; Value = "-19.51%"; IsNumericText = $true },
    @{ Cell = "G47"; Value = "13"; IsNumericText = $true },
    @{ Cell = "D48"; Value = "0.002425"; IsNumericText = $true },
    @{ Cell = "E48"; Value = "13.91%"; IsNumericText = $true },
    @{ Cell = "G48"; Value = "13"; IsNumericText = $true },
    @{ Cell = "E49"; Value = "-0.01%"; IsNumericText = $true },
    @{ Cell = "G49"; Value = "13"; IsNumericText = $true },
    @{ Cell = "E50"; Value = "-0.01%"; IsNumericText = $true },
    @{ Cell = "G50"; Value = "13"; IsNumericText = $true },
    @{ Cell = "G51"; Value = "13"; IsNumericText = $true }
)

foreach ($ch in $changes) {
    $r = $ws.Range($ch.Cell)
    if ($ch.IsNumericText) {
        $r.NumberFormat = "@"
        $r.Value = $ch.Value
        $r.ClearFormats()
    } else {
        $r.Value = $ch.Value
    }
}
